# Update rent sensitivity analysis parameters and adjust rent levels range
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated model parameters (previously: allowance=410, mortgage/expense=2160)
$mgmtFeeRate = 0.15     # Management Fee (£) = Rent * 15%          (unchanged)
$allowance   = 382      # deducted from (Rent - Mgmt Fee) to get Taxable Profit
$taxRate     = 0.4      # Income Tax (£) = Taxable Profit * 40%     (unchanged)
$niRate      = 0.08     # NI (£) = Taxable Profit * 8%              (unchanged)
$taxOffset   = 266      # Total Tax Liability = Income Tax + NI - 266 (unchanged)
$otherCosts  = 2132     # Net Monthly Income = Rent - Mgmt Fee - Total Tax Liability - otherCosts

# Recompute rows 2-15 (rent levels 2000 .. 2650, step 50) with the new parameters
$startRent = 2000
$step = 50
for ($i = 0; $i -lt 14; $i++) {
    $row = $i + 2
    $rent = $startRent + ($i * $step)

    $mgmtFee = $rent * $mgmtFeeRate
    $taxableProfit = $rent - $mgmtFee - $allowance
    $incomeTax = $taxableProfit * $taxRate
    $ni = $taxableProfit * $niRate
    $totalTax = $incomeTax + $ni - $taxOffset
    $netIncome = $rent - $mgmtFee - $totalTax - $otherCosts

    $ws.Cells.Item($row, 3).Value = $taxableProfit
    $ws.Cells.Item($row, 4).Value = $incomeTax
    $ws.Cells.Item($row, 5).Value = $ni
    $ws.Cells.Item($row, 6).Value = $totalTax
    $ws.Cells.Item($row, 7).Value = $netIncome
}

# The rent-level range now stops at 2650 (row 15); remove the old extra rows (16-21)
$ws.Range("A16:G21").Delete()
